$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 1).Value = "2017.05.30 02.59.50"
$ws.Cells.Item(19, 2).Value = 15.14999961853027
$ws.Cells.Item(19, 3).Value = 100
$ws.Cells.Item(19, 4).Value = 1013
$ws.Cells.Item(19, 5).Value = 1.5
$ws.Cells.Item(19, 6).Value = "2017-05-30T09:00:00"
$ws.Cells.Item(19, 7).Value = "2017-05-30T12:00:00"
$ws.Cells.Item(19, 8).Value = 15.07999992370605
$ws.Cells.Item(19, 9).Value = 994.3699951171875
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 0.05999999865889549
$ws.Cells.Item(19, 12).Value = 1.610000014305115

# Row 20
$ws.Cells.Item(20, 1).Value = "2017.05.30 03.00.27"
$ws.Cells.Item(20, 2).Value = 15.14999961853027
$ws.Cells.Item(20, 3).Value = 100
$ws.Cells.Item(20, 4).Value = 1013
$ws.Cells.Item(20, 5).Value = 1.5
$ws.Cells.Item(20, 6).Value = "2017-05-30T09:00:00"
$ws.Cells.Item(20, 7).Value = "2017-05-30T12:00:00"
$ws.Cells.Item(20, 8).Value = 15.07999992370605
$ws.Cells.Item(20, 9).Value = 994.3699951171875
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 0.05999999865889549
$ws.Cells.Item(20, 12).Value = 1.610000014305115
$ws.Cells.Item(20, 13).Value = 16.98
$ws.Cells.Item(20, 14).Value = 2.33
$ws.Cells.Item(20, 15).Value = 84.2

# Row 21
$ws.Cells.Item(21, 1).Value = "2017.05.30 03.00.48"
$ws.Cells.Item(21, 2).Value = 15.14999961853027
$ws.Cells.Item(21, 3).Value = 100
$ws.Cells.Item(21, 4).Value = 1013
$ws.Cells.Item(21, 5).Value = 1.5
$ws.Cells.Item(21, 6).Value = "2017-05-30T09:00:00"
$ws.Cells.Item(21, 7).Value = "2017-05-30T12:00:00"
$ws.Cells.Item(21, 8).Value = 15.07999992370605
$ws.Cells.Item(21, 9).Value = 994.3699951171875
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 0.05999999865889549
$ws.Cells.Item(21, 12).Value = 1.610000014305115
$ws.Cells.Item(21, 13).Value = 17
$ws.Cells.Item(21, 14).Value = 2.34
$ws.Cells.Item(21, 15).Value = 84
